# Update marksheet totals on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct-answer total marks value
$ws.Range("B11").Value = 5

# Total row: aggregate marks obtained
$ws.Range("B12").Value = 105

# Total row: "correct/total" display string
$ws.Range("E12").Value = "105/140"
